$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.192.58"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.43"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.257"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.25"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.878.64"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.643.69"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.538"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.76"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.168.12"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.42"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.43"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.36"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.80"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.58"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.272.27"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.542"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.843"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  +6.49%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.789.61"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.45"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("E48").Value = "  +16.69%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0975"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.03%  "
